# "eksik student mailleri eklendi" -- add the missing "email" column
# (header + per-student mailto: hyperlinks) to the students sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header -----------------------------------------------------------
$ws.Range("F1").Value = "email"

# Column F needs to be wide enough to show the addresses (~19.66 chars).
$ws.Range("F1").ColumnWidth = 18.86

# --- per-row e-mail addresses + hyperlinks -----------------------------
$emails = @{
    2  = "a.yilmaz@live.com"
    3  = "a.demir@live.com"
    4  = "a.korkmaz@live.com"
    5  = "z.durgun@live.com"
    6  = "m.gulmez@live.com"
    7  = "t.guler@live.com"
    8  = "f.agaoglu@live.com"
    9  = "s.guclu@live.com"
    10 = "h.demirev@live.com"
    11 = "m.trabzon@live.com"
}

# row 2
$cell = $ws.Range("F2")
$cell.Value = $emails[2]
$ws.Hyperlinks.Add($cell, "mailto:" + $emails[2], "", "", $emails[2]) | Out-Null

# row 3
$cell = $ws.Range("F3")
$cell.Value = $emails[3]
$ws.Hyperlinks.Add($cell, "mailto:" + $emails[3], "", "", $emails[3]) | Out-Null

# row 4 was first filled/hyperlinked across F4:F11 with the row-3 address
# (a leftover range hyperlink from a fill-down), then corrected per-cell.
$rng = $ws.Range("F4:F11")
$ws.Range("F4").Value = $emails[3]
$ws.Hyperlinks.Add($rng, "mailto:" + $emails[3], "", "", $emails[3]) | Out-Null

# rows 4-11 get their real, individual address + hyperlink
foreach ($r in 4..11) {
    $cell = $ws.Range("F" + $r)
    $cell.Value = $emails[$r]
    $ws.Hyperlinks.Add($cell, "mailto:" + $emails[$r], "", "", $emails[$r]) | Out-Null
}

# --- final selection, matching the saved UI state -----------------------
$ws.Range("F12").Select() | Out-Null
